# Commit: "render website, remove theme (not needed) from docs"
#
# Adds two new paragraph styles to the style sheet ("Abstract Title" and
# "Footnote Block Text") and tightens the space-before on the existing
# "Abstract" style (300 -> 100 twips), matching the pandoc reference.docx
# styles used to render the site.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) New style: "Abstract Title" (w:styleId="AbstractTitle")
#    basedOn Normal, next -> Abstract, centered, bold, small, blue text.
# ---------------------------------------------------------------------
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = "Normal"
$abstractTitle.NextParagraphStyle = "Abstract"
$abstractTitle.QuickStyle = $true

$atPf = $abstractTitle.ParagraphFormat
$atPf.KeepWithNext = $true
$atPf.KeepTogether = $true
$atPf.Alignment = 1            # wdAlignParagraphCenter
$atPf.SpaceBefore = 15         # 300 twips
$atPf.SpaceAfter = 0           # 0 twips

$atFont = $abstractTitle.Font
$atFont.Size = 10              # sz=20 (half-points)
$atFont.SizeBi = 10            # szCs=20
$atFont.Bold = $true
$atFont.Color = 9067060        # 0x345A8A (BGR-packed OLE color)

# ---------------------------------------------------------------------
# 2) Existing style "Abstract": reduce space-before from 300 -> 100
#    (space-after stays at 300)
# ---------------------------------------------------------------------
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5   # 100 twips

# ---------------------------------------------------------------------
# 3) New style: "Footnote Block Text" (w:styleId="FootnoteBlockText")
#    basedOn Footnote Text, next -> Footnote Text, indented block quote
#    layout for footnotes (mirrors the built-in "Block Text" style).
# ---------------------------------------------------------------------
$fnBlockText = $d.Styles.Add("Footnote Block Text", 1)
$fnBlockText.BaseStyle = "Footnote Text"
$fnBlockText.NextParagraphStyle = "Footnote Text"
$fnBlockText.Priority = 9
$fnBlockText.UnhideWhenUsed = $true
$fnBlockText.QuickStyle = $true

$fbtPf = $fnBlockText.ParagraphFormat
$fbtPf.SpaceBefore = 5         # 100 twips
$fbtPf.SpaceAfter = 5          # 100 twips
$fbtPf.FirstLineIndent = 0
$fbtPf.LeftIndent = 24         # 480 twips
$fbtPf.RightIndent = 24        # 480 twips

Write-Output "Styles updated. Total styles now:"
Write-Output $d.Styles.Count
